$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

# Insert the two new rows at the bottom of the table. Inserting (rather than
# just typing past the end of the used range) makes Excel carry the
# column/row formatting down from the row above, same as a user pressing
# Ctrl+D / dragging the fill handle would.
$ws.Rows.Item(151).Insert()
$ws.Rows.Item(152).Insert()

# --- Row 151: 2023-11-20, no bet placed that day ---
$ws.Range("A151").Value = 150
$ws.Range("C151").Value = 45250
$ws.Range("D151").Value = "2023-11-20"

# --- Row 152: 2023-11-21, new bet placed ---
$ws.Range("A152").Value = 151
$ws.Range("B152").Value = 143
$ws.Range("C152").Value = 45251
$ws.Range("D152").Value = "2023-11-21"
$ws.Range("E152").Value = 1
$ws.Range("F152").Value = 1.1
$ws.Range("G152").Formula = "=I150"
$ws.Range("H152").Value = 500
$ws.Range("I152").Formula = "=G152+H152"
$ws.Range("J152").Value = "FUTBOL"
$ws.Range("K152").Value = "CLASIFICATORIAS 2026"
$ws.Range("L152").Formula = "=ROUND((I152/`$G`$31-1)*100, 3)+`$L`$29"
